$wb = $excel.ActiveWorkbook

# Sheet1: move selection from A1:C3 to E5
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E5").Select()

# Sheet2: remove the middle "r" column (B), which shifts the old C/D columns
# left to become the new B/C columns; then update the selection to D6.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns.Item(2).Delete()
$ws2.Range("D6").Select()

# Sheet3: add a selection at C14 (previously no <selection> element existed)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C14").Select()

# Restore Sheet3 as the active/visible tab (matches activeTab=2, tabSelected on Sheet3)
$ws3.Activate()
